# Actualización desde MV -datos-
# Adds a new quarterly observation column (BH, "Agosto.2021") and a new
# trailing row (75, series date "01-04-2021") to the quarterly revisions
# table, mirroring how each new vintage column is appended in this
# dataset: every existing series carries its last known value forward
# into the new column, except the most recent series (row 74) which gets
# revised, and the brand-new series (row 75) which only has one data point
# so far.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. New header cell BH1 = "Agosto.2021" --------------------------
# Clone the header style (bold font + border) from the preceding header
# cell (BG1) first, then overwrite the text - this keeps the new header
# cell's formatting identical to its neighbours without inventing a new
# cell style entry.
$ws.Range("BG1").Copy($ws.Range("BH1"))
$ws.Range("BH1").Value = "Agosto.2021"

# ---- 2. Carry every existing row's last value into the new column ----
# Column BG holds the latest known value for every series published so
# far (rows 2-74). The new vintage column BH simply repeats that same
# number for every series that isn't being revised this time around.
$ws.Range("BG2:BG74").Copy($ws.Range("BH2:BH74"))

# ---- 3. Row 74 ("01-01-2021") is revised in this vintage --------------
$ws.Range("BH74").Value = 12265

# ---- 4. Brand-new row 75 ("01-04-2021") -------------------------------
# Use a scratch cell to produce the literal text "01-04-2021" (rather
# than a value Excel would auto-parse as a date) and Copy it into A75 so
# the destination ends up as a plain shared-string cell with no special
# number format, matching the style-less text cells used throughout
# column A.
$scratch = $ws.Cells.Item(1000, 1000)
$scratch.Value = "=""01-04-2021"""
$scratch.Copy($ws.Range("A75"))
$scratch.Clear()

$ws.Range("BH75").Value = 12671
